$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 772, shifting existing rows 772:828 down to 773:829
$ws.Rows.Item(772).Insert()

# Populate the new row 772 with the new weekly price record
$ws.Range("A772").Value = 3
$ws.Range("B772").Value = "Femacal de La Calera"
$ws.Range("C772").Value = "Coquimbo"
$ws.Range("D772").Value = 45265
$ws.Range("E772").Value = 5
$ws.Range("F772").Value = 100112021
$ws.Range("G772").Value = "Ají"
$ws.Range("H772").Value = "Inferno"
$ws.Range("I772").Value = "Primera"
$ws.Range("J772").Value = 25
$ws.Range("K772").Value = 45000
$ws.Range("L772").Value = 45000
$ws.Range("M772").Value = 45000
$ws.Range("N772").Value = "`$/caja 15 kilos"
$ws.Range("O772").Value = "Provincia de Quillota"
$ws.Range("P772").Value = 3000
$ws.Range("Q772").Value = 15
$ws.Range("R772").Value = "Hortaliza"
